# Auto-generated Excel COM-interop script
# Applies scheduled market-price / profit recalculation updates to the
# Famfrit_Profits workbook across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 599.5
$ws.Cells.Item(4, 9).Value = 599.5
$ws.Cells.Item(4, 11).Value = 599.5
$ws.Cells.Item(4, 13).Value = -485.5

$ws.Cells.Item(21, 8).Value = 650
$ws.Cells.Item(21, 9).Value = 733.3333
$ws.Cells.Item(21, 11).Value = 733.3333
$ws.Cells.Item(21, 13).Value = -265.3333

$ws.Cells.Item(23, 8).Value = 650
$ws.Cells.Item(23, 9).Value = 733.3333
$ws.Cells.Item(23, 11).Value = 733.3333
$ws.Cells.Item(23, 13).Value = -499.3333

$ws.Cells.Item(100, 8).Value = 1748.5883
$ws.Cells.Item(100, 9).Value = 942.1
$ws.Cells.Item(100, 11).Value = 942.1
$ws.Cells.Item(100, 13).Value = -401.1

$ws.Cells.Item(108, 8).Value = 48197.168
$ws.Cells.Item(108, 10).Value = 49045.75
$ws.Cells.Item(108, 12).Value = 49045.75
$ws.Cells.Item(108, 14).Value = -56725.75

$ws.Cells.Item(113, 8).Value = 7400.75
$ws.Cells.Item(113, 10).Value = 9732.666999999999
$ws.Cells.Item(113, 12).Value = 9732.666999999999
$ws.Cells.Item(113, 14).Value = -16240.667

$ws.Cells.Item(137, 8).Value = 2462.8
$ws.Cells.Item(137, 9).Value = 2553.6667
$ws.Cells.Item(137, 10).Value = 2326.5
$ws.Cells.Item(137, 11).Value = 7661.000100000001
$ws.Cells.Item(137, 12).Value = 6979.5
$ws.Cells.Item(137, 13).Value = -5111.000100000001
$ws.Cells.Item(137, 14).Value = -12079.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 23263692
$ws.Cells.Item(32, 9).Value = 29416478
$ws.Cells.Item(32, 10).Value = 19833.334
$ws.Cells.Item(32, 11).Value = 29416478
$ws.Cells.Item(32, 12).Value = 19833.334
$ws.Cells.Item(32, 13).Value = -29416191
$ws.Cells.Item(32, 14).Value = -20407.334

$ws.Cells.Item(54, 8).Value = 36747.5
$ws.Cells.Item(54, 10).Value = 36747.5
$ws.Cells.Item(54, 12).Value = 36747.5
$ws.Cells.Item(54, 14).Value = -38285.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2283.25
$ws.Cells.Item(20, 9).Value = 2224.875
$ws.Cells.Item(20, 10).Value = 2400
$ws.Cells.Item(20, 11).Value = 2224.875
$ws.Cells.Item(20, 12).Value = 2400
$ws.Cells.Item(20, 13).Value = -1977.875
$ws.Cells.Item(20, 14).Value = -2894

$ws.Cells.Item(134, 8).Value = 3005.4
$ws.Cells.Item(134, 9).Value = 2873.4324
$ws.Cells.Item(134, 11).Value = 8620.297200000001
$ws.Cells.Item(134, 13).Value = -6085.297200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 354.3684
$ws.Cells.Item(7, 9).Value = 298.7143
$ws.Cells.Item(7, 10).Value = 386.83334
$ws.Cells.Item(7, 11).Value = 298.7143
$ws.Cells.Item(7, 12).Value = 386.83334
$ws.Cells.Item(7, 13).Value = -185.7143
$ws.Cells.Item(7, 14).Value = -612.83334

$ws.Cells.Item(22, 8).Value = 11627.556
$ws.Cells.Item(22, 9).Value = 20299.8
$ws.Cells.Item(22, 10).Value = 787.25
$ws.Cells.Item(22, 11).Value = 20299.8
$ws.Cells.Item(22, 12).Value = 787.25
$ws.Cells.Item(22, 13).Value = -19949.8
$ws.Cells.Item(22, 14).Value = -1487.25

$ws.Cells.Item(31, 8).Value = 24394672
$ws.Cells.Item(31, 9).Value = 3338.7307
$ws.Cells.Item(31, 11).Value = 3338.7307
$ws.Cells.Item(31, 13).Value = -3043.7307

$ws.Cells.Item(34, 8).Value = 24394672
$ws.Cells.Item(34, 9).Value = 3338.7307
$ws.Cells.Item(34, 11).Value = 3338.7307
$ws.Cells.Item(34, 13).Value = -3136.7307

$ws.Cells.Item(108, 8).Value = 22861
$ws.Cells.Item(108, 9).Value = 23905.25
$ws.Cells.Item(108, 10).Value = 18684
$ws.Cells.Item(108, 11).Value = 23905.25
$ws.Cells.Item(108, 12).Value = 18684
$ws.Cells.Item(108, 13).Value = -20065.25
$ws.Cells.Item(108, 14).Value = -26364

$ws.Cells.Item(134, 8).Value = 1439.5
$ws.Cells.Item(134, 9).Value = 1216.7858
$ws.Cells.Item(134, 11).Value = 3650.3574
$ws.Cells.Item(134, 13).Value = -1115.3574

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 30703238
$ws.Cells.Item(4, 9).Value = 36280996
$ws.Cells.Item(4, 11).Value = 108842988
$ws.Cells.Item(4, 13).Value = -108842876

$ws.Cells.Item(56, 8).Value = 23770.1
$ws.Cells.Item(56, 9).Value = 23770.1
$ws.Cells.Item(56, 11).Value = 23770.1
$ws.Cells.Item(56, 13).Value = -23240.1

$ws.Cells.Item(88, 8).Value = 19507
$ws.Cells.Item(88, 10).Value = 19507
$ws.Cells.Item(88, 12).Value = 58521
$ws.Cells.Item(88, 14).Value = -59377

$ws.Cells.Item(91, 8).Value = 19507
$ws.Cells.Item(91, 10).Value = 19507
$ws.Cells.Item(91, 12).Value = 58521
$ws.Cells.Item(91, 14).Value = -61485

$ws.Cells.Item(96, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 2025
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 12).Value = 0

$ws.Cells.Item(101, 8).Value = 7981
$ws.Cells.Item(101, 10).Value = 7981
$ws.Cells.Item(101, 12).Value = 23943
$ws.Cells.Item(101, 14).Value = -28811

$ws.Cells.Item(103, 8).Value = 2032.6666
$ws.Cells.Item(103, 10).Value = 2998.5
$ws.Cells.Item(103, 12).Value = 8995.5
$ws.Cells.Item(103, 14).Value = -10753.5

$ws.Cells.Item(128, 8).Value = 115987.75
$ws.Cells.Item(128, 9).Value = 115987.75
$ws.Cells.Item(128, 11).Value = 347963.25
$ws.Cells.Item(128, 13).Value = -342983.25

$ws.Cells.Item(131, 8).Value = 1747.3939
$ws.Cells.Item(131, 9).Value = 1365.4546
$ws.Cells.Item(131, 10).Value = 1938.3636
$ws.Cells.Item(131, 11).Value = 4096.3638
$ws.Cells.Item(131, 12).Value = 5815.0908
$ws.Cells.Item(131, 13).Value = 943.6361999999999
$ws.Cells.Item(131, 14).Value = -15895.0908

$ws.Cells.Item(133, 8).Value = 10477.958
$ws.Cells.Item(133, 9).Value = 5311.273
$ws.Cells.Item(133, 10).Value = 14849.77
$ws.Cells.Item(133, 11).Value = 15933.819
$ws.Cells.Item(133, 12).Value = 44549.31
$ws.Cells.Item(133, 13).Value = -10873.819
$ws.Cells.Item(133, 14).Value = -54669.31

$ws.Cells.Item(134, 8).Value = 4349.316
$ws.Cells.Item(134, 9).Value = 1702.8667
$ws.Cells.Item(134, 10).Value = 14273.5
$ws.Cells.Item(134, 11).Value = 5108.6001
$ws.Cells.Item(134, 12).Value = 42820.5
$ws.Cells.Item(134, 13).Value = -38.60009999999966
$ws.Cells.Item(134, 14).Value = -52960.5

$ws.Cells.Item(137, 8).Value = 6105.385
$ws.Cells.Item(137, 10).Value = 7647.75
$ws.Cells.Item(137, 12).Value = 22943.25
$ws.Cells.Item(137, 14).Value = -33143.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(69, 14).ClearContents()
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0

$ws.Cells.Item(72, 14).ClearContents()
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3362.25
$ws.Cells.Item(40, 9).Value = 3634.2
$ws.Cells.Item(40, 11).Value = 3634.2
$ws.Cells.Item(40, 13).Value = -3498.2

$ws.Cells.Item(46, 8).Value = 1409.3405
$ws.Cells.Item(46, 9).Value = 675.5454999999999
$ws.Cells.Item(46, 11).Value = 675.5454999999999
$ws.Cells.Item(46, 13).Value = -487.5454999999999

$ws.Cells.Item(82, 8).Value = 6584.7144
$ws.Cells.Item(82, 9).Value = 4049.5
$ws.Cells.Item(82, 10).Value = 7598.8
$ws.Cells.Item(82, 11).Value = 4049.5
$ws.Cells.Item(82, 12).Value = 7598.8
$ws.Cells.Item(82, 13).Value = -3688.5
$ws.Cells.Item(82, 14).Value = -8320.799999999999

$ws.Cells.Item(85, 8).Value = 6584.7144
$ws.Cells.Item(85, 9).Value = 4049.5
$ws.Cells.Item(85, 10).Value = 7598.8
$ws.Cells.Item(85, 11).Value = 4049.5
$ws.Cells.Item(85, 12).Value = 7598.8
$ws.Cells.Item(85, 13).Value = -2801.5
$ws.Cells.Item(85, 14).Value = -10094.8

$ws.Cells.Item(122, 8).Value = 4385.657
$ws.Cells.Item(122, 9).Value = 4065.261
$ws.Cells.Item(122, 10).Value = 4999.75
$ws.Cells.Item(122, 11).Value = 12195.783
$ws.Cells.Item(122, 12).Value = 14999.25
$ws.Cells.Item(122, 13).Value = -9745.782999999999
$ws.Cells.Item(122, 14).Value = -19899.25

$ws.Cells.Item(136, 8).Value = 2279.9636
$ws.Cells.Item(136, 9).Value = 1715.96
$ws.Cells.Item(136, 10).Value = 7920
$ws.Cells.Item(136, 11).Value = 5147.88
$ws.Cells.Item(136, 12).Value = 23760
$ws.Cells.Item(136, 13).Value = -2597.88
$ws.Cells.Item(136, 14).Value = -28860

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 8).Value = 18123.75
$ws.Cells.Item(40, 10).Value = 21498.334
$ws.Cells.Item(40, 12).Value = 21498.334
$ws.Cells.Item(40, 14).Value = -21796.334

$ws.Cells.Item(113, 8).Value = 616.25
$ws.Cells.Item(113, 9).Value = 241.11111
$ws.Cells.Item(113, 11).Value = 723.3333299999999
$ws.Cells.Item(113, 13).Value = 1446.66667

$ws.Cells.Item(126, 8).Value = 3535.0322
$ws.Cells.Item(126, 9).Value = 3912.4443
$ws.Cells.Item(126, 11).Value = 11737.3329
$ws.Cells.Item(126, 13).Value = -9267.332900000001

$ws.Cells.Item(132, 8).Value = 4118.795
$ws.Cells.Item(132, 9).Value = 4201.1313
$ws.Cells.Item(132, 11).Value = 12603.3939
$ws.Cells.Item(132, 13).Value = -10073.3939
